$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-15"

# Update the header label in I1 ("2022 (through 06-14)" -> "2022 (through 06-15)")
$ws.Range("I1").Value = "2022 (through 06-15)"

# Update June 2022 value (row 7 = June)
$ws.Range("I7").Value = 69

# Update November 2021 value (row 12 = November)
$ws.Range("H12").Value = 201

# Update Totals row (row 14)
$ws.Range("H14").Value = 1848
$ws.Range("I14").Value = 732
